$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column D (Duracion) from time-fraction numbers (h:mm format) to plain text "m:ss" strings
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D2:D118").NumberFormat = "@"

$durations = @{
    12 = "2:45"
    13 = "3:47"
    14 = "3:31"
    15 = "2:50"
    16 = "3:54"
    17 = "3:44"
    18 = "2:45"
    19 = "3:20"
    20 = "3:17"
    21 = "2:36"
    22 = "2:55"
    23 = "3:17"
    24 = "5:30"
    25 = "3:31"
    26 = "2:40"
    27 = "3:50"
    28 = "3:59"
    29 = "3:43"
    30 = "2:49"
    31 = "7:10"
    32 = "4:08"
    33 = "3:35"
    34 = "4:20"
    35 = "3:59"
    36 = "5:36"
    37 = "3:10"
    38 = "3:32"
    39 = "3:27"
    40 = "3:40"
    41 = "3:17"
    42 = "3:14"
    43 = "3:06"
    44 = "4:10"
    45 = "4:00"
    46 = "3:45"
    47 = "2:50"
    48 = "4:45"
    49 = "4:01"
    50 = "4:17"
    51 = "3:25"
    52 = "2:48"
    53 = "3:14"
    54 = "3:19"
    55 = "3:04"
    56 = "3:21"
    57 = "2:34"
    58 = "2:54"
    59 = "3:19"
    60 = "3:16"
    61 = "3:11"
    62 = "2:50"
    63 = "2:41"
    64 = "2:59"
    65 = "2:57"
    66 = "2:30"
    67 = "2:35"
    68 = "2:39"
    69 = "2:39"
    70 = "2:42"
    71 = "2:56"
    72 = "3:10"
    73 = "2:26"
    74 = "2:30"
    75 = "3:08"
    76 = "2:41"
    77 = "3:10"
    78 = "3:05"
    79 = "2:43"
    88 = "10:24"
    89 = "6:41"
    90 = "6:00"
    91 = "3:44"
    92 = "3:10"
    93 = "4:52"
    94 = "2:31"
    95 = "3:49"
    96 = "5:23"
    97 = "4:19"
    98 = "4:11"
    99 = "4:27"
    100 = "2:33"
    101 = "3:45"
    102 = "4:00"
    103 = "3:47"
    104 = "5:15"
    105 = "2:55"
    106 = "8:18"
    107 = "6:13"
    108 = "4:01"
    109 = "3:11"
    110 = "5:15"
}
foreach ($r in $durations.Keys) {
    $ws.Range("D$r").Value = $durations[$r]
}

# Add/update "Notas" (column N) content
$notes = @{
    2 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    3 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    4 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    5 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    6 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    7 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    8 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    9 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    10 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    11 = "Producido por Joe Cain / Arreglos: Javier Vasquez"
    12 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    13 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    14 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    15 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    16 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    17 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    18 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    19 = "Producido por  Louie Ramirez/ Coro: Adalberto Santiago, Ruben Blades, Yayo el Indio y Carlos Ayala"
    20 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    21 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    22 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    23 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    24 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    25 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    26 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    27 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    28 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    29 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    30 = "Producido por Joe Cain, Miguel Estivill, Tito Puente / Arreglos: Javier Vasquez"
    42 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    43 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    44 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    45 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    46 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    47 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    48 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    49 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    50 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    51 = "Producido por Joe Cain, Luis Ramirez, Ismael Rivera, Pancho Cristal"
    52 = "Producido por Pancho Cristal"
    53 = "Producido por Pancho Cristal"
    54 = "Producido por Pancho Cristal"
    55 = "Producido por Pancho Cristal"
    56 = "Producido por Pancho Cristal"
    57 = "Producido por Pancho Cristal"
    58 = "Producido por Pancho Cristal"
    59 = "Producido por Pancho Cristal"
    60 = "Producido por Pancho Cristal"
    61 = "Producido por Pancho Cristal"
    62 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    63 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    64 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    65 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    66 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    67 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    68 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    69 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    70 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    71 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    72 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    73 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    74 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    75 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    76 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    77 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    78 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    79 = "Con: Ismael Rivera, Rafael Ithier, Roberto Roena, Sammy Ayala, Kito Velez, Hector Santos…"
    87 = "Producido por Willie Colón"
    103 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    104 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    105 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    106 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    107 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    108 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    109 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
    110 = "Producido por Willie Colón & Jerry Masucci/ Arreglos: Wllie Colon/ Coro: Justo Betancourt, Johnny Pacheco"
}
foreach ($r in $notes.Keys) {
    $ws.Range("N$r").Value = $notes[$r]
}

# Update column N width to fit the longer note text
$ws.Columns("N:N").ColumnWidth = 91

# Restore view/selection state
$ws.Range("N97").Select()